$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197. This shifts the existing rows 197-239
# down to 198-240 (old row 239 becomes row 240), matching the target diff.
$ws.Rows("197:197").Insert()

# Populate the newly inserted row 197 with the new record's data.
$ws.Range("A197").Value = 5
$ws.Range("B197").Value = "Macroferia Regional de Talca"
$ws.Range("C197").Value = "Maule"
$ws.Range("D197").Value = 45015
$ws.Range("E197").Value = 7
$ws.Range("F197").Value = 100112031
$ws.Range("G197").Value = "Poroto verde"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 200
$ws.Range("K197").Value = 30000
$ws.Range("L197").Value = 30000
$ws.Range("M197").Value = 30000
$ws.Range("N197").Value = "$/saco 25 kilos"
$ws.Range("O197").Value = "Región Metropolitana"
$ws.Range("P197").Value = 1200
$ws.Range("Q197").Value = 25
$ws.Range("R197").Value = "Hortaliza"
